$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 1545329.1
$ws.Range("I62").Value2 = 2576876
$ws.Range("J62").Value2 = 169933.33
$ws.Range("K62").Value2 = 2576876
$ws.Range("L62").Value2 = 169933.33
$ws.Range("M62").Value2 = -2576252
$ws.Range("N62").Value2 = -171181.33
$ws.Range("H65").Value2 = 1545329.1
$ws.Range("I65").Value2 = 2576876
$ws.Range("J65").Value2 = 169933.33
$ws.Range("K65").Value2 = 12884380
$ws.Range("L65").Value2 = 849666.6499999999
$ws.Range("M65").Value2 = -12881260
$ws.Range("N65").Value2 = -855906.6499999999
$ws.Range("H76").Value2 = 6605109
$ws.Range("I76").Value2 = 58698.223
$ws.Range("K76").Value2 = 58698.223
$ws.Range("M76").Value2 = -58383.223
$ws.Range("H79").Value2 = 6605109
$ws.Range("I79").Value2 = 58698.223
$ws.Range("K79").Value2 = 58698.223
$ws.Range("M79").Value2 = -57606.223
$ws.Range("H88").Value2 = 2049.75
$ws.Range("I88").Value2 = 1849.5
$ws.Range("K88").Value2 = 1849.5
$ws.Range("M88").Value2 = -1443.5
$ws.Range("H91").Value2 = 2049.75
$ws.Range("I91").Value2 = 1849.5
$ws.Range("K91").Value2 = 1849.5
$ws.Range("M91").Value2 = -445.5
$ws.Range("H106").Value2 = 151485.72
$ws.Range("J106").Value2 = 341166.66
$ws.Range("L106").Value2 = 341166.66
$ws.Range("N106").Value2 = -342428.66
$ws.Range("H133").Value2 = 0
$ws.Range("J133").Value2 = 0
$ws.Range("L133").Value2 = 0
$ws.Range("N133").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4140.5835
$ws.Range("I32").Value2 = 4154.241
$ws.Range("J32").Value2 = 4084
$ws.Range("K32").Value2 = 4154.241
$ws.Range("L32").Value2 = 4084
$ws.Range("M32").Value2 = -3867.241
$ws.Range("N32").Value2 = -4658
$ws.Range("H61").Value2 = 30306992
$ws.Range("I61").Value2 = 47621420
$ws.Range("K61").Value2 = 47621420
$ws.Range("M61").Value2 = -47621208
$ws.Range("H132").Value2 = 90911690
$ws.Range("I132").Value2 = 111113620
$ws.Range("K132").Value2 = 333340860
$ws.Range("M132").Value2 = -333338330
$ws.Range("H136").Value2 = 30306992
$ws.Range("I136").Value2 = 47621420
$ws.Range("K136").Value2 = 142864260
$ws.Range("M136").Value2 = -142861710

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 54999
$ws.Range("I20").Value2 = 99998
$ws.Range("J20").Value2 = 10000
$ws.Range("K20").Value2 = 99998
$ws.Range("L20").Value2 = 10000
$ws.Range("M20").Value2 = -99751
$ws.Range("N20").Value2 = -10494
$ws.Range("H86").Value2 = 3809.375
$ws.Range("I86").Value2 = 4720.5557
$ws.Range("J86").Value2 = 2637.8572
$ws.Range("K86").Value2 = 4720.5557
$ws.Range("L86").Value2 = 2637.8572
$ws.Range("M86").Value2 = -3597.5557
$ws.Range("N86").Value2 = -4883.8572
$ws.Range("H89").Value2 = 3809.375
$ws.Range("I89").Value2 = 4720.5557
$ws.Range("J89").Value2 = 2637.8572
$ws.Range("K89").Value2 = 23602.7785
$ws.Range("L89").Value2 = 13189.286
$ws.Range("M89").Value2 = -17986.7785
$ws.Range("N89").Value2 = -24421.286
$ws.Range("H94").Value2 = 3625.4167
$ws.Range("I94").Value2 = 3000.625
$ws.Range("K94").Value2 = 3000.625
$ws.Range("M94").Value2 = -2549.625
$ws.Range("H105").Value2 = 2036.1818
$ws.Range("I105").Value2 = 1955.3334
$ws.Range("K105").Value2 = 1955.3334
$ws.Range("M105").Value2 = -208.3334
$ws.Range("H134").Value2 = 6946856
$ws.Range("I134").Value2 = 8335561.5
$ws.Range("K134").Value2 = 25006684.5
$ws.Range("M134").Value2 = -25004149.5

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2963.8472
$ws.Range("I31").Value2 = 2191.9697
$ws.Range("J31").Value2 = 3453.6924
$ws.Range("K31").Value2 = 2191.9697
$ws.Range("L31").Value2 = 3453.6924
$ws.Range("M31").Value2 = -1896.9697
$ws.Range("N31").Value2 = -4043.6924
$ws.Range("H34").Value2 = 2963.8472
$ws.Range("I34").Value2 = 2191.9697
$ws.Range("J34").Value2 = 3453.6924
$ws.Range("K34").Value2 = 2191.9697
$ws.Range("L34").Value2 = 3453.6924
$ws.Range("M34").Value2 = -1989.9697
$ws.Range("N34").Value2 = -3857.6924
$ws.Range("H107").Value2 = 2427.6428
$ws.Range("I107").Value2 = 1629.5217
$ws.Range("K107").Value2 = 1629.5217
$ws.Range("M107").Value2 = 290.4783
$ws.Range("H122").Value2 = 2929.353
$ws.Range("I122").Value2 = 2527.7334
$ws.Range("J122").Value2 = 5941.5
$ws.Range("K122").Value2 = 7583.2002
$ws.Range("L122").Value2 = 17824.5
$ws.Range("M122").Value2 = -5133.2002
$ws.Range("N122").Value2 = -22724.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 1372.8572
$ws.Range("I5").Value2 = 458
$ws.Range("K5").Value2 = 1374
$ws.Range("M5").Value2 = -1262
$ws.Range("H12").Value2 = 328.125
$ws.Range("J12").Value2 = 381
$ws.Range("L12").Value2 = 1143
$ws.Range("N12").Value2 = -1489
$ws.Range("H107").Value2 = 723.93335
$ws.Range("J107").Value2 = 489.2857
$ws.Range("L107").Value2 = 1467.8571
$ws.Range("N107").Value2 = -5307.8571
$ws.Range("H113").Value2 = 771.0476
$ws.Range("J113").Value2 = 1015.53845
$ws.Range("L113").Value2 = 3046.61535
$ws.Range("N113").Value2 = -7386.61535
$ws.Range("H129").Value2 = 2134.8
$ws.Range("I129").Value2 = 560
$ws.Range("J129").Value2 = 4497
$ws.Range("K129").Value2 = 1680
$ws.Range("L129").Value2 = 13491
$ws.Range("M129").Value2 = 3320
$ws.Range("N129").Value2 = -23491
$ws.Range("H130").Value2 = 999
$ws.Range("I130").Value2 = 999
$ws.Range("K130").Value2 = 2997
$ws.Range("M130").Value2 = 2023
$ws.Range("H132").Value2 = 1399.3334
$ws.Range("J132").Value2 = 1500
$ws.Range("L132").Value2 = 13500
$ws.Range("N132").Value2 = -18560
$ws.Range("H135").Value2 = 1372.8572
$ws.Range("I135").Value2 = 458
$ws.Range("K135").Value2 = 4122
$ws.Range("M135").Value2 = -1587

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 1874.1177
$ws.Range("I80").Value2 = 2029.5555
$ws.Range("J80").Value2 = 1699.25
$ws.Range("K80").Value2 = 2029.5555
$ws.Range("L80").Value2 = 1699.25
$ws.Range("M80").Value2 = -1031.5555
$ws.Range("N80").Value2 = -3695.25
$ws.Range("H83").Value2 = 1874.1177
$ws.Range("I83").Value2 = 2029.5555
$ws.Range("J83").Value2 = 1699.25
$ws.Range("K83").Value2 = 10147.7775
$ws.Range("L83").Value2 = 8496.25
$ws.Range("M83").Value2 = -5155.7775
$ws.Range("N83").Value2 = -18480.25
$ws.Range("H107").Value2 = 1713.8572
$ws.Range("I107").Value2 = 499.42856
$ws.Range("J107").Value2 = 2928.2856
$ws.Range("K107").Value2 = 499.42856
$ws.Range("L107").Value2 = 2928.2856
$ws.Range("M107").Value2 = 1420.57144
$ws.Range("N107").Value2 = -6768.2856
$ws.Range("H132").Value2 = 3559.7856
$ws.Range("I132").Value2 = 2955.5
$ws.Range("J132").Value2 = 4365.5
$ws.Range("K132").Value2 = 8866.5
$ws.Range("L132").Value2 = 13096.5
$ws.Range("M132").Value2 = -6336.5
$ws.Range("N132").Value2 = -18156.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 997.25
$ws.Range("I46").Value2 = 997.25
$ws.Range("J46").Value2 = 0
$ws.Range("K46").Value2 = 997.25
$ws.Range("L46").Value2 = 0
$ws.Range("M46").Value2 = -809.25
$ws.Range("N46").ClearContents()
$ws.Range("H132").Value2 = 2588.8655
$ws.Range("I132").Value2 = 2389.3022
$ws.Range("K132").Value2 = 7167.9066
$ws.Range("M132").Value2 = -4637.9066
$ws.Range("H136").Value2 = 2557.8696
$ws.Range("I136").Value2 = 2084.0625
$ws.Range("J136").Value2 = 3640.8572
$ws.Range("K136").Value2 = 6252.1875
$ws.Range("L136").Value2 = 10922.5716
$ws.Range("M136").Value2 = -3702.1875
$ws.Range("N136").Value2 = -16022.5716
$ws.Range("H139").Value2 = 91000
$ws.Range("I139").Value2 = 91000
$ws.Range("K139").Value2 = 91000
$ws.Range("M139").Value2 = -85860

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value2 = 21250
$ws.Range("I26").Value2 = 17500
$ws.Range("K26").Value2 = 17500
$ws.Range("M26").Value2 = -17207
$ws.Range("H100").Value2 = 1089.0834
$ws.Range("I100").Value2 = 952.7143
$ws.Range("K100").Value2 = 1905.4286
$ws.Range("M100").Value2 = -1364.4286
$ws.Range("H126").Value2 = 1679.3334
$ws.Range("I126").Value2 = 1632.8889
$ws.Range("J126").Value2 = 1749
$ws.Range("K126").Value2 = 4898.6667
$ws.Range("L126").Value2 = 5247
$ws.Range("M126").Value2 = -2428.6667
$ws.Range("N126").Value2 = -10187
$ws.Range("H132").Value2 = 1718.8948
$ws.Range("I132").Value2 = 1486.0857
$ws.Range("J132").Value2 = 4435
$ws.Range("K132").Value2 = 4458.257100000001
$ws.Range("L132").Value2 = 13305
$ws.Range("M132").Value2 = -1928.257100000001
$ws.Range("N132").Value2 = -18365
$ws.Range("H136").Value2 = 3107.8215
$ws.Range("I136").Value2 = 1557.2106
$ws.Range("J136").Value2 = 6381.3335
$ws.Range("K136").Value2 = 4671.6318
$ws.Range("L136").Value2 = 19144.0005
$ws.Range("M136").Value2 = -2121.6318
$ws.Range("N136").Value2 = -24244.0005
$ws.Range("H137").Value2 = 59999
$ws.Range("J137").Value2 = 59999
$ws.Range("L137").Value2 = 59999
$ws.Range("N137").Value2 = -70199
$ws.Range("H138").Value2 = 94000
$ws.Range("J138").Value2 = 94000
$ws.Range("L138").Value2 = 94000
$ws.Range("N138").Value2 = -104280
